$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.348.77"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "'1.840.22"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'239.00"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'0.6261"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07431"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "'0.2895"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'24.84"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").Value = "'0.07712"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'1.839.11"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "'4.951"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'0.6743"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'0.00001022"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "'81.70"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "'6.214"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'29.340.04"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'231.65"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'7.365"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'158.15"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "'8.462"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'0.1343"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "'0.07280"
$ws.Range("E28").Value = "  +12.81%  "
$ws.Range("D30").Value = "'1.476"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'4.039"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "'4.040"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").Value = "'1.818"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "'0.6955"
$ws.Range("D36").Value = "'2.567"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").Value = "'6.920"
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("D38").Value = "'0.01834"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'1.232.58"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").Value = "'0.9462"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'1.984.51"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "'100.56"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "'65.40"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("E46").Value = "  +4.81%  "
$ws.Range("D47").Value = "'1.713"
$ws.Range("E47").Value = "  -2.48%  "
$ws.Range("D48").Value = "'6.939"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "'8.907"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'0.1136"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").Value = "'0.3901"
$ws.Range("E51").Value = "  -1.19%  "
